$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Finish grading the "CustomerMapping Class" rubric table (rows 18-25):
# award full credit (10 points) for the "whoPurchasedProduct() method" row (22)
# and the "findAllBrands()" row (24) in the "Total Points" column E.
$ws.Range("E22").Value = 10
$ws.Range("E24").Value = 10

# Leave the sheet scrolled/selected where the grader finished working.
$ws.Range("E24").Select()
